$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fndc5"
$ws.Range("C2").Value = "Itgb5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05565899999999999
$ws.Range("H2").Value = 0.166977
$ws.Range("I2").Value = 0.01431489546586528
$ws.Range("J2").Value = 0.01431489546586528
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.460162333333333
$ws.Range("N2").Value = 28.380487
$ws.Range("O2").Value = 0.08197024919772995
$ws.Range("P2").Value = 0.08197024919772995
$ws.Range("Q2").Value = 0.5265431753109999
$ws.Range("R2").Value = 4.738888577798999
$ws.Range("S2").Value = 0.001173395548576432
$ws.Range("T2").Value = 0.001173395548576432

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fndc5"
$ws.Range("C3").Value = "Itgb5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05565899999999999
$ws.Range("H3").Value = 0.166977
$ws.Range("I3").Value = 0.01431489546586528
$ws.Range("J3").Value = 0.01431489546586528
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 46.18256633333333
$ws.Range("N3").Value = 138.547699
$ws.Range("O3").Value = 0.4001618933742075
$ws.Range("P3").Value = 0.4001618933742075
$ws.Range("Q3").Value = 2.570475459547
$ws.Range("R3").Value = 23.134279135923
$ws.Range("S3").Value = 0.005728275673074511
$ws.Range("T3").Value = 0.005728275673074511

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fndc5"
$ws.Range("C4").Value = "Itgb5"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05565899999999999
$ws.Range("H4").Value = 0.166977
$ws.Range("I4").Value = 0.01431489546586528
$ws.Range("J4").Value = 0.01431489546586528
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 51.40166966666666
$ws.Range("N4").Value = 154.205009
$ws.Range("O4").Value = 0.4453842886934318
$ws.Range("P4").Value = 0.4453842886934319
$ws.Range("Q4").Value = 2.860965531977
$ws.Range("R4").Value = 25.748689787793
$ws.Range("S4").Value = 0.006375629534785242
$ws.Range("T4").Value = 0.006375629534785243

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fndc5"
$ws.Range("C5").Value = "Itgb5"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05565899999999999
$ws.Range("H5").Value = 0.166977
$ws.Range("I5").Value = 0.01431489546586528
$ws.Range("J5").Value = 0.01431489546586528
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.365307333333334
$ws.Range("N5").Value = 25.095922
$ws.Range("O5").Value = 0.07248356873463072
$ws.Range("P5").Value = 0.07248356873463073
$ws.Range("Q5").Value = 0.465604640866
$ws.Range("R5").Value = 4.190441767794
$ws.Range("S5").Value = 0.0010375947094291
$ws.Range("T5").Value = 0.0010375947094291

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fndc5"
$ws.Range("C6").Value = "Itgb5"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.907196
$ws.Range("H6").Value = 5.721588000000001
$ws.Range("I6").Value = 0.4905102745812252
$ws.Range("J6").Value = 0.4905102745812252
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.460162333333333
$ws.Range("N6").Value = 28.380487
$ws.Range("O6").Value = 0.08197024919772995
$ws.Range("P6").Value = 0.08197024919772995
$ws.Range("Q6").Value = 18.042383761484
$ws.Range("R6").Value = 162.381453853356
$ws.Range("S6").Value = 0.04020724944146997
$ws.Range("T6").Value = 0.04020724944146997

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fndc5"
$ws.Range("C7").Value = "Itgb5"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.907196
$ws.Range("H7").Value = 5.721588000000001
$ws.Range("I7").Value = 0.4905102745812252
$ws.Range("J7").Value = 0.4905102745812252
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 46.18256633333333
$ws.Range("N7").Value = 138.547699
$ws.Range("O7").Value = 0.4001618933742075
$ws.Range("P7").Value = 0.4001618933742075
$ws.Range("Q7").Value = 88.079205780668
$ws.Range("R7").Value = 792.712852026012
$ws.Range("S7").Value = 0.1962835201959255
$ws.Range("T7").Value = 0.1962835201959255

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fndc5"
$ws.Range("C8").Value = "Itgb5"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.907196
$ws.Range("H8").Value = 5.721588000000001
$ws.Range("I8").Value = 0.4905102745812252
$ws.Range("J8").Value = 0.4905102745812252
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 51.40166966666666
$ws.Range("N8").Value = 154.205009
$ws.Range("O8").Value = 0.4453842886934318
$ws.Range("P8").Value = 0.4453842886934319
$ws.Range("Q8").Value = 98.033058781588
$ws.Range("R8").Value = 882.297529034292
$ws.Range("S8").Value = 0.2184655697411789
$ws.Range("T8").Value = 0.2184655697411789

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fndc5"
$ws.Range("C9").Value = "Itgb5"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.907196
$ws.Range("H9").Value = 5.721588000000001
$ws.Range("I9").Value = 0.4905102745812252
$ws.Range("J9").Value = 0.4905102745812252
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 8.365307333333334
$ws.Range("N9").Value = 25.095922
$ws.Range("O9").Value = 0.07248356873463072
$ws.Range("P9").Value = 0.07248356873463073
$ws.Range("Q9").Value = 15.954280684904
$ws.Range("R9").Value = 143.588526164136
$ws.Range("S9").Value = 0.03555393520265082
$ws.Range("T9").Value = 0.03555393520265083

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Fndc5"
$ws.Range("C10").Value = "Itgb5"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3444803333333333
$ws.Range("H10").Value = 1.033441
$ws.Range("I10").Value = 0.08859663238134165
$ws.Range("J10").Value = 0.08859663238134165
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.460162333333333
$ws.Range("N10").Value = 28.380487
$ws.Range("O10").Value = 0.08197024919772995
$ws.Range("P10").Value = 0.08197024919772995
$ws.Range("Q10").Value = 3.25883987397411
$ws.Range("R10").Value = 29.32955886576699
$ws.Range("S10").Value = 0.007262288034378245
$ws.Range("T10").Value = 0.007262288034378245

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Fndc5"
$ws.Range("C11").Value = "Itgb5"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3444803333333333
$ws.Range("H11").Value = 1.033441
$ws.Range("I11").Value = 0.08859663238134165
$ws.Range("J11").Value = 0.08859663238134165
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 46.18256633333333
$ws.Range("N11").Value = 138.547699
$ws.Range("O11").Value = 0.4001618933742075
$ws.Range("P11").Value = 0.4001618933742075
$ws.Range("Q11").Value = 15.90898584469544
$ws.Range("R11").Value = 143.180872602259
$ws.Range("S11").Value = 0.0354529961602963
$ws.Range("T11").Value = 0.0354529961602963

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Fndc5"
$ws.Range("C12").Value = "Itgb5"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3444803333333333
$ws.Range("H12").Value = 1.033441
$ws.Range("I12").Value = 0.08859663238134165
$ws.Range("J12").Value = 0.08859663238134165
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 51.40166966666666
$ws.Range("N12").Value = 154.205009
$ws.Range("O12").Value = 0.4453842886934318
$ws.Range("P12").Value = 0.4453842886934319
$ws.Range("Q12").Value = 17.70686430066322
$ws.Range("R12").Value = 159.361778705969
$ws.Range("S12").Value = 0.03945954809379732
$ws.Range("T12").Value = 0.03945954809379732

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Fndc5"
$ws.Range("C13").Value = "Itgb5"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3444803333333333
$ws.Range("H13").Value = 1.033441
$ws.Range("I13").Value = 0.08859663238134165
$ws.Range("J13").Value = 0.08859663238134165
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 8.365307333333334
$ws.Range("N13").Value = 25.095922
$ws.Range("O13").Value = 0.07248356873463072
$ws.Range("P13").Value = 0.07248356873463073
$ws.Range("Q13").Value = 2.881683858622444
$ws.Range("R13").Value = 25.935154727602
$ws.Range("S13").Value = 0.006421800092869787
$ws.Range("T13").Value = 0.006421800092869788

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Fndc5"
$ws.Range("C14").Value = "Itgb5"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.580852333333333
$ws.Range("H14").Value = 4.742557
$ws.Range("I14").Value = 0.4065781975715678
$ws.Range("J14").Value = 0.4065781975715678
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 9.460162333333333
$ws.Range("N14").Value = 28.380487
$ws.Range("O14").Value = 0.08197024919772995
$ws.Range("P14").Value = 0.08197024919772995
$ws.Range("Q14").Value = 14.95511969836211
$ws.Range("R14").Value = 134.596077285259
$ws.Range("S14").Value = 0.03332731617330529
$ws.Range("T14").Value = 0.03332731617330529

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Fndc5"
$ws.Range("C15").Value = "Itgb5"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.580852333333333
$ws.Range("H15").Value = 4.742557
$ws.Range("I15").Value = 0.4065781975715678
$ws.Range("J15").Value = 0.4065781975715678
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 46.18256633333333
$ws.Range("N15").Value = 138.547699
$ws.Range("O15").Value = 0.4001618933742075
$ws.Range("P15").Value = 0.4001618933742075
$ws.Range("Q15").Value = 73.00781774737145
$ws.Range("R15").Value = 657.070359726343
$ws.Range("S15").Value = 0.1626971013449112
$ws.Range("T15").Value = 0.1626971013449112

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Fndc5"
$ws.Range("C16").Value = "Itgb5"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.580852333333333
$ws.Range("H16").Value = 4.742557
$ws.Range("I16").Value = 0.4065781975715678
$ws.Range("J16").Value = 0.4065781975715678
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 51.40166966666666
$ws.Range("N16").Value = 154.205009
$ws.Range("O16").Value = 0.4453842886934318
$ws.Range("P16").Value = 0.4453842886934319
$ws.Range("Q16").Value = 81.25844942977922
$ws.Range("R16").Value = 731.3260448680129
$ws.Range("S16").Value = 0.1810835413236703
$ws.Range("T16").Value = 0.1810835413236703

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Fndc5"
$ws.Range("C17").Value = "Itgb5"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.580852333333333
$ws.Range("H17").Value = 4.742557
$ws.Range("I17").Value = 0.4065781975715678
$ws.Range("J17").Value = 0.4065781975715678
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 8.365307333333334
$ws.Range("N17").Value = 25.095922
$ws.Range("O17").Value = 0.07248356873463072
$ws.Range("P17").Value = 0.07248356873463073
$ws.Range("Q17").Value = 13.22431561695045
$ws.Range("R17").Value = 119.018840552554
$ws.Range("S17").Value = 0.029470238729681
$ws.Range("T17").Value = 0.02947023872968101

